$d = $word.ActiveDocument

function Insert-SuperscriptLabel {
    param(
        [int]$footnoteIndex,
        [string]$searchText,
        [string]$label,
        [int]$colorBgr
    )
    $fn = $d.Footnotes.Item($footnoteIndex)
    $searchFrom = $fn.Reference.End
    $docEnd = $d.Content.End
    $r = $d.Range($searchFrom, $docEnd)
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $startPos = $r.Start
    $r.Text = $label + $searchText
    $labelRange = $d.Range($startPos, $startPos + $label.Length)
    $labelRange.Font.Color = $colorBgr
    $labelRange.Font.Size = 10
    $labelRange.Font.Superscript = $true
}

# Colors (VBA BGR-packed decimal):
#   E74C3C -> 3951847
#   3498DB -> 14391348
#   9B59B6 -> 11950491

Insert-SuperscriptLabel 1 "。師者，所以傳道受業解惑也" "[注]" 3951847
Insert-SuperscriptLabel 2 "。人非生而知之者，孰能無惑" "[解讀]" 14391348
Insert-SuperscriptLabel 3 "？惑而不從師，其為惑也" "[注]" 3951847
Insert-SuperscriptLabel 4 "，終不解矣" "[解讀]" 14391348
Insert-SuperscriptLabel 5 "。" "[脂批]" 11950491

# Remove numbering from footnote labels (keep only the base label)
$d.Footnotes.Item(1).Range.Text = " 【注】師：老師，指有專門知識或技能的人。"
$d.Footnotes.Item(2).Range.Text = " 【解讀】這句話說明老師的三個職責：傳授道理、教授學業、解答疑惑。"
$d.Footnotes.Item(3).Range.Text = " 【注】孰：誰。"
$d.Footnotes.Item(4).Range.Text = " 【解讀】有疑惑卻不向老師請教，那疑惑就永遠無法解開。"
$d.Footnotes.Item(5).Range.Text = " 【脂批】此開卷第一回也。作者自云曾歷過一番夢幻之後，故將真事隱去。"
